$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells L1:N1 (copy formatting from the existing K1 header) ---
$ws.Range("K1").Copy($ws.Range("L1:N1"))
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# --- New data columns L, M, N for rows 2-7 ---
$lValues = @{ 2 = 92.07932629902824; 3 = 86.28489520037526; 4 = 88.79569566345432; 5 = 100.9425353902489; 6 = 19.4629156234702; 7 = 21.90935307361503 }
$mValues = @{ 2 = 228260; 3 = 35293; 4 = 180657; 5 = 22989; 6 = 2121; 7 = 87 }
$nValues = @{ 2 = 296.8270481144343; 3 = 578.5737704918033; 4 = 141.2486317435497; 5 = 221.0480769230769; 6 = 14.14; 7 = 43.5 }

for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 12).Value = $lValues[$r]
    $ws.Cells.Item($r, 13).Value = $mValues[$r]
    $ws.Cells.Item($r, 14).Value = $nValues[$r]
}

# --- Rescale existing "particip" (E) and "taxa_sucesso" (F) columns from
#     fractional (0-1) values to percentage-point values (0-100) ---
$eValues = @{ 2 = 94.9812734082397; 3 = 5.018726591760299; 4 = 92.91553133514986; 5 = 7.084468664850137; 6 = 99.70760233918129; 7 = 0.2923976608187134 }
$fValues = @{ 2 = 60.64668769716089; 3 = 91.04477611940298; 4 = 93.76832844574781; 5 = 100; 6 = 21.9941348973607; 7 = 100 }

for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 5).Value = $eValues[$r]
    $ws.Cells.Item($r, 6).Value = $fValues[$r]
}

Write-Host ("FINAL E2 value2: " + $ws.Cells.Item(2,5).Value2)
Write-Host ("FINAL L2 value2: " + $ws.Cells.Item(2,12).Value2)
